# Refresh the cryptos price/volume table (Price = column D, Volume(1h) = column E).
# Values are stored as text (e.g. "39.887.18", "  +0.26%  "), so every write is
# apostrophe-prefixed to stop Excel from reinterpreting numeric-looking text as
# a real number, and the style is reset to "Normal" afterwards so the
# quote-prefix flag doesn't linger as a visible cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.887.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.26%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.208.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.79%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'288.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'86.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.64%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.34%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.20%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'30.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.46%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0775"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.57%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.45%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.70%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.549.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'13.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.62%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.214.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.42%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.52%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'39.820.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.28%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'11.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +10.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.47%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.01%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'65.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'235.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.87%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E27").Value = "'  -1.92%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.85%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'152.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.91%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'31.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.97%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.10%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.00%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.73%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.81%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'15.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.30%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.36%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +2.57%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.41%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.085.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.95%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.37%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +6.19%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.20%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'17.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +7.05%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.19%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.422.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.78%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'88.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.75%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.57%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'68.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.98%  "
$ws.Range("E51").Style = "Normal"
